$wb = $excel.ActiveWorkbook

# --- Sheet "Encabezado" ---
$ws1 = $wb.Worksheets.Item("Encabezado")

$ws1.Range("B1").Value = "01135036"
$ws1.Range("B2").Value = "T010005584"
$ws1.Range("B3").Value = "0001"
$ws1.Range("B4").Value = "T3C1"
$ws1.Range("B5").Value = "VSUAREZ"
$ws1.Range("B6").Value = "CT3000000000001"
$ws1.Range("B7").Value = "T03"
$ws1.Range("B8").Value = 36.5635
$ws1.Range("B10").Value = "0.00"
$ws1.Range("B11").Value = "0.00"

# --- Sheet "Detalles" ---
$ws2 = $wb.Worksheets.Item("Detalles")

# Update row 2 values
$ws2.Range("B2").Value = "03000038"
$ws2.Range("C2").Value = "1.00"
$ws2.Range("E2").Value = "AGUACATE  010075"

# Remove rows 3 through 6 (data previously there is no longer needed)
$ws2.Range("A3:I6").Delete()
